$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the new shared strings in the exact order they appear in the target
# file (index 18..23), so the resulting sharedStrings.xml table lines up
# with the diff: total delta-v=, MF, burn dv, finmass, avg acce, MJ time.
$ws.Range("G6").Value = "total delta-v = "
$ws.Range("G2").Value = "MF"
$ws.Range("G9").Value = "burn dv"
$ws.Range("G13").Value = "finmass"
$ws.Range("G15").Value = "avg acce"
$ws.Range("G17").Value = "MJ time"

# Numeric input and formulas
$ws.Range("G3").Value = 0.45
$ws.Range("G7").Formula = "=E1*E2*LN(E3/G3)"
$ws.Range("G10").Formula = "=SQRT(B2/B15)*(1-SQRT(1-B14))"
$ws.Range("H12").Formula = "=G10/G7"
$ws.Range("H13").Formula = "=E3/EXP(LN(E3/G3)*H12)"
$ws.Range("H15").Formula = "=E4/((E3+H13)/2)"
$ws.Range("H17").Formula = "=G10/H15"

# Update selection to match target state
$ws.Range("G18").Select()
